$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row was added: A12 gets a single space value.
$ws.Range("A12").Value = " "

# Selection moved to C1.
$ws.Range("C1").Select() | Out-Null

# Column A width changed slightly (e.g. via auto-fit / resize).
$ws.Columns("A:A").ColumnWidth = 82.7

# Page setup explicitly touched (paper size + portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "Done"
